$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.167.17"
$ws.Range("D3").Value = "1.899.98"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.27"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5247"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E8").Value = "  +0.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07290"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.36"
$ws.Range("E10").Value = "  +1.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9043"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08147"
$ws.Range("E12").Value = "  -3.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.45"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.348"
$ws.Range("D15").Value = "1.794.83"
$ws.Range("E15").Value = "  -5.71%  "
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.70"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "27.205.29"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.104"
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.80"
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.327"
$ws.Range("E24").Value = "  +2.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.38"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.21"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.828"
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.877"
$ws.Range("E30").Value = "  -0.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09226"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05062"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7941"
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.224"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.978"
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.368"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.655"
$ws.Range("E37").Value = "  +2.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5706"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01989"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  +0.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.007"
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.575"
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.17"
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4876"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.629"
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.52"
$ws.Range("E49").Value = "  +2.86%  "
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("E51").Value = "  +0.47%  "
